{"js": "// Replace the multiplication-problem text in each table cell with its new\n// value, per the diff. Each \"old\" string is unique in the document, so a\n// body-wide search-and-replace (matchCase to avoid any ambiguity) is safe.\nconst replacements = [\n  [\"357\u00d76=\", \"827\u00d76=\"],\n  [\"689\u00d72=\", \"934\u00d79=\"],\n  [\"810\u00d75=\", \"129\u00d73=\"],\n  [\"407\u00d75=\", \"265\u00d74=\"],\n  [\"987\u00d78=\", \"739\u00d78=\"],\n  [\"623\u00d72=\", \"866\u00d78=\"],\n  [\"590\u00d79=\", \"960\u00d79=\"],\n  [\"549\u00d76=\", \"744\u00d76=\"],\n  [\"752\u00d72=\", \"975\u00d78=\"],\n  [\"886\u00d73=\", \"225\u00d73=\"],\n  [\"609\u00d76=\", \"678\u00d79=\"],\n  [\"314\u00d77=\", \"527\u00d79=\"],\n  [\"470\u00d78=\", \"913\u00d76=\"],\n  [\"250\u00d77=\", \"898\u00d76=\"],\n  [\"978\u00d76=\", \"867\u00d75=\"],\n  [\"771\u00d77=\", \"893\u00d74=\"],\n  [\"201\u00d77=\", \"987\u00d75=\"],\n  [\"209\u00d78=\", \"208\u00d78=\"],\n  [\"675\u00d78=\", \"552\u00d74=\"],\n  [\"861\u00d72=\", \"171\u00d77=\"],\n  [\"131\u00d74=\", \"781\u00d72=\"],\n  [\"403\u00d72=\", \"447\u00d74=\"],\n  [\"589\u00d78=\", \"900\u00d78=\"],\n  [\"942\u00d76=\", \"538\u00d72=\"],\n  [\"901\u00d76=\", \"795\u00d72=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the multiplication-problem text in each table cell with its new\n# value, per the diff. Each \"old\" string is unique in the document, so a\n# simple Find/Replace (ReplaceAll, MatchCase on) for each pair is safe.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"357\u00d76=\", \"827\u00d76=\"),\n    @(\"689\u00d72=\", \"934\u00d79=\"),\n    @(\"810\u00d75=\", \"129\u00d73=\"),\n    @(\"407\u00d75=\", \"265\u00d74=\"),\n    @(\"987\u00d78=\", \"739\u00d78=\"),\n    @(\"623\u00d72=\", \"866\u00d78=\"),\n    @(\"590\u00d79=\", \"960\u00d79=\"),\n    @(\"549\u00d76=\", \"744\u00d76=\"),\n    @(\"752\u00d72=\", \"975\u00d78=\"),\n    @(\"886\u00d73=\", \"225\u00d73=\"),\n    @(\"609\u00d76=\", \"678\u00d79=\"),\n    @(\"314\u00d77=\", \"527\u00d79=\"),\n    @(\"470\u00d78=\", \"913\u00d76=\"),\n    @(\"250\u00d77=\", \"898\u00d76=\"),\n    @(\"978\u00d76=\", \"867\u00d75=\"),\n    @(\"771\u00d77=\", \"893\u00d74=\"),\n    @(\"201\u00d77=\", \"987\u00d75=\"),\n    @(\"209\u00d78=\", \"208\u00d78=\"),\n    @(\"675\u00d78=\", \"552\u00d74=\"),\n    @(\"861\u00d72=\", \"171\u00d77=\"),\n    @(\"131\u00d74=\", \"781\u00d72=\"),\n    @(\"403\u00d72=\", \"447\u00d74=\"),\n    @(\"589\u00d78=\", \"900\u00d78=\"),\n    @(\"942\u00d76=\", \"538\u00d72=\"),\n    @(\"901\u00d76=\", \"795\u00d72=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
